$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-11 23:48:50'
$ws.Range('H2').NumberFormat = "@"
$ws.Range('H2').Value = '84%'
$ws.Range('I2').Value = '9.1 mm'
$ws.Range('O2').Value = '2.6 °C'
$ws.Range('E3').Value = '2026-02-11 23:48:53'
$ws.Range('H3').NumberFormat = "@"
$ws.Range('H3').Value = '79%'
$ws.Range('I3').Value = '4.7 mm'
$ws.Range('E4').Value = '2026-02-11 23:48:56'
$ws.Range('J4').Value = '1001.5 hPa'
$ws.Range('E5').Value = '2026-02-11 23:48:59'
$ws.Range('I5').Value = '3.9 mm'
$ws.Range('E6').Value = '2026-02-11 23:49:01'
$ws.Range('J6').Value = '1002.0 hPa'
$ws.Range('E7').Value = '2026-02-11 23:49:04'
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H7').Value = '49%'
$ws.Range('J7').Value = '1002.9 hPa'
$ws.Range('O7').Value = '18.4 °C'
$ws.Range('E8').Value = '2026-02-11 23:49:07'
$ws.Range('J8').Value = '1002.6 hPa'
$ws.Range('O8').Value = '14.6 °C'
$ws.Range('E9').Value = '2026-02-11 23:49:09'
$ws.Range('H9').NumberFormat = "@"
$ws.Range('H9').Value = '90%'
$ws.Range('E10').Value = '2026-02-11 23:49:12'
$ws.Range('H10').NumberFormat = "@"
$ws.Range('H10').Value = '78%'
$ws.Range('O10').Value = '13.5 °C'
$ws.Range('E11').Value = '2026-02-11 23:49:14'
$ws.Range('O11').Value = '7.8 °C'
$ws.Range('E12').Value = '2026-02-11 23:49:17'
$ws.Range('E13').Value = '2026-02-11 23:49:20'
$ws.Range('I13').Value = '2.5 mm'
$ws.Range('J13').Value = '1004.5 hPa'
$ws.Range('E14').Value = '2026-02-11 23:49:23'
$ws.Range('H14').NumberFormat = "@"
$ws.Range('H14').Value = '56%'
$ws.Range('N14').Value = '9.2 °C 23:29 TU'
$ws.Range('O14').Value = '17.9 °C'
$ws.Range('E15').Value = '2026-02-11 23:49:25'
$ws.Range('E16').Value = '2026-02-11 23:49:28'
$ws.Range('I16').Value = '7.9 mm'
$ws.Range('E17').Value = '2026-02-11 23:49:31'
$ws.Range('H17').NumberFormat = "@"
$ws.Range('H17').Value = '81%'
$ws.Range('E18').Value = '2026-02-11 23:49:34'
$ws.Range('J18').Value = '1002.1 hPa'
$ws.Range('E19').Value = '2026-02-11 23:49:36'
$ws.Range('E20').Value = '2026-02-11 23:49:39'
$ws.Range('E21').Value = '2026-02-11 23:49:42'
$ws.Range('J21').Value = '1004.9 hPa'
$ws.Range('E22').Value = '2026-02-11 23:49:44'
$ws.Range('I22').Value = '6.0 mm'
$ws.Range('E23').Value = '2026-02-11 23:49:47'
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H23').Value = '77%'
$ws.Range('I23').Value = '7.5 mm'
$ws.Range('E24').Value = '2026-02-11 23:49:49'
$ws.Range('J24').Value = '1006.0 hPa'
$ws.Range('N24').Value = '9.9 °C 23:29 TU'
$ws.Range('O24').Value = '12.8 °C'
$ws.Range('E25').Value = '2026-02-11 23:49:52'
$ws.Range('H25').NumberFormat = "@"
$ws.Range('H25').Value = '69%'
$ws.Range('N25').Value = '-1.1 °C 23:29 TU'
$ws.Range('E26').Value = '2026-02-11 23:49:55'
$ws.Range('J26').Value = '1001.9 hPa'
$ws.Range('E27').Value = '2026-02-11 23:49:58'
$ws.Range('E28').Value = '2026-02-11 23:50:01'
$ws.Range('J28').Value = '1002.2 hPa'
$ws.Range('E29').Value = '2026-02-11 23:50:03'
$ws.Range('E30').Value = '2026-02-11 23:50:06'
$ws.Range('J30').Value = '1002.1 hPa'
$ws.Range('E31').Value = '2026-02-11 23:50:08'
$ws.Range('J31').Value = '1001.3 hPa'
$ws.Range('O31').Value = '14.3 °C'
$ws.Range('E32').Value = '2026-02-11 23:50:11'
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H32').Value = '79%'
$ws.Range('E33').Value = '2026-02-11 23:50:14'
$ws.Range('J33').Value = '1004.1 hPa'
$ws.Range('O33').Value = '6.4 °C'
$ws.Range('E34').Value = '2026-02-11 23:50:16'
$ws.Range('E35').Value = '2026-02-11 23:50:19'
$ws.Range('E36').Value = '2026-02-11 23:50:22'
$ws.Range('J36').Value = '1002.2 hPa'
$ws.Range('O36').Value = '13.1 °C'
$ws.Range('E37').Value = '2026-02-11 23:50:24'
$ws.Range('J37').Value = '1003.6 hPa'
$ws.Range('E38').Value = '2026-02-11 23:50:27'
$ws.Range('I38').Value = '4.2 mm'
$ws.Range('O38').Value = '15.3 °C'
$ws.Range('E39').Value = '2026-02-11 23:50:30'
$ws.Range('O39').Value = '1.2 °C'
$ws.Range('E40').Value = '2026-02-11 23:50:33'
$ws.Range('J40').Value = '1006.0 hPa'
$ws.Range('O40').Value = '7.5 °C'
$ws.Range('E41').Value = '2026-02-11 23:50:35'
$ws.Range('E42').Value = '2026-02-11 23:50:38'
$ws.Range('E43').Value = '2026-02-11 23:50:41'
$ws.Range('E44').Value = '2026-02-11 23:50:43'
$ws.Range('I44').Value = '7.8 mm'
$ws.Range('E45').Value = '2026-02-11 23:50:46'
$ws.Range('H45').NumberFormat = "@"
$ws.Range('H45').Value = '88%'
$ws.Range('J45').Value = '1004.6 hPa'
$ws.Range('E46').Value = '2026-02-11 23:50:48'
$ws.Range('H46').NumberFormat = "@"
$ws.Range('H46').Value = '64%'
$ws.Range('J46').Value = '1006.3 hPa'
$ws.Range('L46').Value = '54.7 km/h - 268º 23:00 TU'
